$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "worn wear"
$ws.Range("A2").Value = "wow clothes for women"
$ws.Range("A3").Value = "write tight"
$ws.Range("A4").Value = "x 3"
$ws.Range("A5").Value = "x compression pants"
$ws.Range("A6").Value = "x compression shorts"
$ws.Range("A7").Value = "x endurance"
$ws.Range("A8").Value = "x factor workout bands"
$ws.Range("A9").Value = "x fit"
$ws.Range("A10").Value = "x gear"
$ws.Range("A11").Value = "x knee brace"
$ws.Range("A12").Value = "x large exercise ball"
$ws.Range("A13").Value = "x marks the spot"
$ws.Range("A14").Value = "x small"
$ws.Range("A15").Value = "x sport"
$ws.Range("A16").Value = "xl back brace"
$ws.Range("A17").Value = "xl basketball hoop"
$ws.Range("A18").Value = "xl compression pants"
$ws.Range("A19").Value = "xl knee braces for women"
$ws.Range("A20").Value = "xl knee support"
$ws.Range("A21").Value = "xl leggings for women"
$ws.Range("A22").Value = "xl winter leggings for women"
$ws.Range("A23").Value = "xl womens winter clothes"
$ws.Range("A24").Value = "xl yoga capris"
$ws.Range("A25").Value = "xlarge leggings for women"
$ws.Range("A26").Value = "xmas tights"
$ws.Range("A27").Value = "xs compression pants"
$ws.Range("A28").Value = "xs knee compression"
$ws.Range("A29").Value = "xs stockings"
$ws.Range("A30").Value = "xs womens bike"
$ws.Range("A31").Value = "xs womens leggings"
$ws.Range("A32").Value = "xsmall leggings for women"
$ws.Range("A33").Value = "xxl workout"
$ws.Range("A34").Value = "xxs clothes women"
$ws.Range("A35").Value = "xxs pants for women"
$ws.Range("A36").Value = "y leggings"
$ws.Range("A37").Value = "yoga after knee replacement"
$ws.Range("A38").Value = "yoga athletic pants"
$ws.Range("A39").Value = "yoga ball base ring"
$ws.Range("A40").Value = "yoga ball with stability ring"
$ws.Range("A41").Value = "yoga capri"
$ws.Range("A42").Value = "yoga capri pants for women"
$ws.Range("A43").Value = "yoga capris for women"
$ws.Range("A44").Value = "yoga capris for women high waist"
$ws.Range("A45").Value = "yoga clotges"
$ws.Range("A46").Value = "yoga clothes"
$ws.Range("A47").Value = "yoga clothes for women"
$ws.Range("A48").Value = "yoga clothing for women"
$ws.Range("A49").Value = "yoga compression leggings"
$ws.Range("A50").Value = "yoga compression shorts for women"
$ws.Range("A51").Value = "yoga conditioning"
$ws.Range("A52").Value = "yoga exercise clothes"
$ws.Range("A53").Value = "yoga fitness"
$ws.Range("A54").Value = "yoga for hip replacement"
$ws.Range("A55").Value = "yoga for runners"
$ws.Range("A56").Value = "yoga gear women"
$ws.Range("A57").Value = "yoga gym"
$ws.Range("A58").Value = "yoga knee support"
$ws.Range("A59").Value = "yoga leggigns"
$ws.Range("A60").Value = "yoga leggings high waist"
$ws.Range("A61").Value = "yoga leggings knee length"
$ws.Range("A62").Value = "yoga leggings over foot"
$ws.Range("A63").Value = "yoga pant capris"
$ws.Range("A64").Value = "yoga pant underwear"
$ws.Range("A65").Value = "yoga pants amazon choice"
$ws.Range("A66").Value = "yoga pants compression"
$ws.Range("A67").Value = "yoga pants drawstring capri"
$ws.Range("A68").Value = "yoga pants extra long length"
$ws.Range("A69").Value = "yoga pants for woman stretch leggings"
$ws.Range("A70").Value = "yoga pants for women capri length"
$ws.Range("A71").Value = "yoga pants for women long length"
$ws.Range("A72").Value = "yoga pants for women stretch leggings fitness running"
$ws.Range("A73").Value = "yoga pants for women stretch leggings fitness running sports active"
$ws.Range("A74").Value = "yoga pants for women victoria secret"
$ws.Range("A75").Value = "yoga pants gym"
$ws.Range("A76").Value = "yoga pants knee length"
$ws.Range("A77").Value = "yoga pants knee length for women"
$ws.Range("A78").Value = "yoga pants long length"
$ws.Range("A79").Value = "yoga pants running pants"
$ws.Range("A80").Value = "yoga pants short length"
$ws.Range("A81").Value = "yoga pants tall length"
$ws.Range("A82").Value = "yoga pants tight"
$ws.Range("A83").Value = "yoga pants tight for women"
$ws.Range("A84").Value = "yoga pants with back pockets"
$ws.Range("A85").Value = "yoga pants women long"
$ws.Range("A86").Value = "yoga pants you can wear to work"
$ws.Range("A87").Value = "yoga pelvis"
$ws.Range("A88").Value = "yoga people"
$ws.Range("A89").Value = "yoga pro plus"
$ws.Range("A90").Value = "yoga rings for women"
$ws.Range("A91").Value = "yoga tights"
$ws.Range("A92").Value = "yoga tights for women"
$ws.Range("A93").Value = "yoga tights high waist"
$ws.Range("A94").Value = "yoga tights with pockets"
$ws.Range("A95").Value = "yoga waist buttery soft leggings"
$ws.Range("A96").Value = "yoga waist leggings"
$ws.Range("A97").Value = "yoga wear"
$ws.Range("A98").Value = "yoga wear women"
$ws.Range("A99").Value = "yoga workout capris"
$ws.Range("A100").Value = "yogi clothes for women"
